$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new user record (row 33) by duplicating the formatting of the
# last existing data row (32) so the new row picks up the same cell
# styles (left-aligned boolean in I, filled style in D), then overwrite
# the values with the new "Ewan Marsh" record.
$ws.Rows("32").Copy()
$ws.Rows("33").Insert(-4121)

$ws.Range("A33").Value = 110032
$ws.Range("B33").Value = 9317596770
$ws.Range("C33").Value = "Ewan Marsh"
$ws.Range("D33").Value = "ewan.marsh@xyz.com"
$ws.Range("E33").Value = 818876433
$ws.Range("F33").Value = "ACT"
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = "PWD"
$ws.Range("I33").Value = $true
$ws.Range("J33").Value = "superadmin"
$ws.Range("K33").Value = "now()"

# Match the author's final selection/scroll state: columns L onward
# selected (e.g. after pressing Ctrl+Space on col L), view scrolled back
# to the top-left.
$ws.Range("L1:XFD1048576").Select()
